# cronograma.xlsx update:
#  - insere nova linha "Estudo piloto" (Set/2011) logo apos "Analise do perfil..."
#  - "Entrevistas" passa a abranger tambem Out/2011 (coluna D), alem de Nov/2011
#  - acrescenta nova linha "Escrita de artigo para journal" (Mar/2012-Abr/2012) ao final
#  - acrescenta coluna de mes "Abr/2012"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Abrir espaco para a nova linha 7 ("Estudo piloto") empurrando as linhas
#    7..14 para 8..15. Faz-se de baixo para cima copiando cada linha inteira
#    (valores + formatacao) para a linha seguinte, o que preserva os indices
#    de estilo ja existentes (ao contrario de Rows.Insert, que cria novos).
# ---------------------------------------------------------------------------
for ($r = 14; $r -ge 7; $r--) {
    $destRow = $r + 1
    $ws.Range("A" + $r + ":I" + $r).Copy($ws.Range("A" + $destRow + ":I" + $destRow))
}

# A linha 7 (ainda com o conteudo antigo da linha 7, "Execucao do experimento")
# e reformatada como a linha 6 (mesmo mes, Set/2011) e recebe o novo texto.
$ws.Range("A6:I6").Copy($ws.Range("A7:I7"))
$ws.Range("A7").Value = "Estudo piloto"

# "Entrevistas" (agora na linha 11) passa a comecar em Out/2011 tambem.
$ws.Range("D11").Style = $ws.Range("E11").Style

# ---------------------------------------------------------------------------
# 2) Nova linha 16 ao final: "Escrita de artigo para journal"
#    (formatacao igual a da linha 15, que ja tem marcadores em H e I).
# ---------------------------------------------------------------------------
$ws.Range("A15:I15").Copy($ws.Range("A16:I16"))
$ws.Range("A16").Value = "Escrita de artigo para journal"
$ws.Range("H16").Style = $ws.Range("H14").Style
$ws.Range("I16").Style = $ws.Range("I15").Style

# ---------------------------------------------------------------------------
# 3) Nova coluna J: mes "Abr/2012"
# ---------------------------------------------------------------------------
$ws.Range("I2:I16").Copy($ws.Range("J2:J16"))
$ws.Range("J2").Value = "Abr/2012"
$ws.Range("J16").Style = $ws.Range("I16").Style

# ---------------------------------------------------------------------------
# 4) Selecao final (a dimensao e recalculada automaticamente ao salvar)
# ---------------------------------------------------------------------------
$ws.Range("D22").Select()
